$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: name / size / iteration
$names = @("US11", "US12", "US13", "US21", "US22", "US23", "US24", "US31")
$sizes = @(1, 3, 5, 1, 3, 5, 8, 3)
$iterations = @("Iteration 1", "Iteration 1", "Iteration 1", "Iteration 2", "Iteration 2", "Iteration 2", "Iteration 2", "Iteration 3")

# Fill column A (name) first
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

# Then column B (size)
for ($i = 0; $i -lt $sizes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $sizes[$i]
}

# Then column C (iteration)
for ($i = 0; $i -lt $iterations.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $iterations[$i]
}

# Move selection to the next empty row in column C, matching where the
# author's cursor landed after finishing data entry.
[void]$ws.Range("C10").Select()
